$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'58.872.55"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -3.40%  "

# Row 3
$ws.Range("D3").Value = "'2.568.10"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.34%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'508.66"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.71%  "

# Row 6
$ws.Range("D6").Value = "'143.81"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -7.09%  "

# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("D8").Value = "'0.556"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -5.96%  "

# Row 9
$ws.Range("D9").Value = "'2.572.98"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.39%  "

# Row 10
$ws.Range("D10").Value = "'6.22"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -7.27%  "

# Row 11
$ws.Range("E11").Value = "  -3.36%  "

# Row 12
$ws.Range("D12").Value = "'0.331"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.56%  "

# Row 13
$ws.Range("E13").Value = "  -0.95%  "

# Row 14
$ws.Range("D14").Value = "'3.013.46"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.50%  "

# Row 15
$ws.Range("D15").Value = "'58.833.82"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.49%  "

# Row 16
$ws.Range("D16").Value = "'20.61"
$ws.Range("D16").ClearFormats()

# Row 17
$ws.Range("E17").Value = "  -4.71%  "

# Row 18
$ws.Range("D18").Value = "'2.569.26"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.34%  "

# Row 19
$ws.Range("E19").Value = "  -5.15%  "

# Row 20
$ws.Range("D20").Value = "'333.56"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -6.04%  "

# Row 21
$ws.Range("D21").Value = "'10.08"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.59%  "

# Row 22
$ws.Range("E22").Value = "  -0.01%  "

# Row 23
$ws.Range("D23").Value = "'5.96"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.26%  "

# Row 24
$ws.Range("D24").Value = "'59.78"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.93%  "

# Row 25
$ws.Range("D25").Value = "'0.407"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.58%  "

# Row 26
$ws.Range("E26").Value = "  +0.15%  "

# Row 27
$ws.Range("E27").Value = "  -5.53%  "

# Row 28
$ws.Range("D28").Value = "'0.0₃0781"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -7.88%  "

# Row 29
$ws.Range("D29").Value = "'6.90"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -7.06%  "

# Row 31
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'149.43"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.89%  "

# Row 32
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "'5.85"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -6.92%  "

# Row 33
$ws.Range("D33").Value = "'18.57"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.32%  "

# Row 34
$ws.Range("E34").Value = "  -3.66%  "

# Row 35
$ws.Range("D35").Value = "'3.94"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.94%  "

# Row 36
$ws.Range("D36").Value = "'0.896"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.14%  "

# Row 37
$ws.Range("E37").Value = "  -8.22%  "

# Row 38
$ws.Range("D38").Value = "'35.94"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.50%  "

# Row 39
$ws.Range("D39").Value = "'0.825"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.83%  "

# Row 40
$ws.Range("D40").Value = "'288.56"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.61%  "

# Row 41
$ws.Range("E41").Value = "  -7.98%  "

# Row 42
$ws.Range("D42").Value = "'3.52"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -7.51%  "

# Row 43
$ws.Range("D43").Value = "'0.997"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.07%  "

# Row 44
$ws.Range("E44").Value = "  -2.39%  "

# Row 45
$ws.Range("E45").Value = "  -3.44%  "

# Row 46
$ws.Range("E46").Value = "  -4.80%  "

# Row 47
$ws.Range("D47").Value = "'18.81"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.84%  "

# Row 48
$ws.Range("D48").Value = "'10.33"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.06%  "

# Row 49
$ws.Range("E49").Value = "  -4.55%  "

# Row 50
$ws.Range("E50").Value = "  -8.08%  "

# Row 51
$ws.Range("D51").Value = "'1.917.17"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.44%  "
